# Apply the "cryptos list" refresh described in the commit diff.
# Each coin row is addressed directly by its cell reference; values that
# look like plain decimal numbers (single "." separator) are prefixed with
# a leading apostrophe so Excel stores them as text (matching the original
# inlineStr cell type) instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = "25.887.28"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "1.640.56"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'213.44"
$ws.Range("E5").Value = "  +2.57%  "

$ws.Range("D6").Value = "'0.5209"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'0.2602"
$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("D9").Value = "'0.06318"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "'20.58"
$ws.Range("E10").Value = "  -1.96%  "

$ws.Range("D11").Value = "'0.07666"
$ws.Range("E11").Value = "  +1.62%  "

$ws.Range("D12").Value = "1.641.12"
$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("D13").Value = "'4.409"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("D14").Value = "1.864.02"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("D15").Value = "'0.5477"
$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").Value = "0.0₅8175"
$ws.Range("E16").Value = "  +2.78%  "

$ws.Range("D17").Value = "'64.42"

$ws.Range("D18").Value = "25.906.40"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").Value = "'4.684"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").Value = "'188.29"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "'10.14"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").Value = "'6.240"
$ws.Range("E23").Value = "  +0.74%  "

$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("D25").Value = "'142.95"
$ws.Range("E25").Value = "  -4.07%  "

$ws.Range("D26").Value = "'0.1236"
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").Value = "'7.355"
$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("D29").Value = "'1.406"
$ws.Range("E29").Value = "  +3.33%  "

$ws.Range("D30").Value = "'0.05927"
$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("D31").Value = "'1.257"
$ws.Range("E31").Value = "  -0.53%  "

$ws.Range("D32").Value = "'3.388"
$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("D33").Value = "'3.397"
$ws.Range("E33").Value = "  -2.05%  "

$ws.Range("D34").Value = "'1.635"
$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").Value = "'0.9862"
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").Value = "'2.396"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D38").Value = "'0.5598"
$ws.Range("E38").Value = "  -4.74%  "

$ws.Range("D39").Value = "'0.01598"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").Value = "'5.814"
$ws.Range("E40").Value = "  -3.07%  "

$ws.Range("D41").Value = "'0.8510"
$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "1.019.86"
$ws.Range("E43").Value = "  -7.90%  "

$ws.Range("D44").Value = "'98.52"
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("D45").Value = "1.788.79"
$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'55.39"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.003"

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.052"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05137"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4208"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'5.875"
$ws.Range("E51").Value = "  +0.09%  "
